# Update countries & provincias Spain
# Applies the 1-May-2020 19:52 data refresh to the "Pais" sheet:
#  - bumps the "last updated" timestamp string
#  - updates several countries' case counters
#  - Guinea's totals overtake Cuba's, so the two swap table rows (80/81)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 19:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1111510
$ws.Range("C4").Value = 16487
$ws.Range("E4").Value = 888817
$ws.Range("G4").Value = 1028
$ws.Range("H4").Value = 64884

# Row 15 - Canada
$ws.Range("B15").Value = 53670
$ws.Range("C15").Value = 434
$ws.Range("D15").Value = 22095
$ws.Range("E15").Value = 28351

# Row 37 - Rumania
$ws.Range("E37").Value = 7495
$ws.Range("G37").Value = 27
$ws.Range("H37").Value = 744

# Rows 80/81 - Guinea overtakes Cuba in total cases, swapping their order
$ws.Range("A80").Value = "Guinea"
$ws.Range("B80").Value = 1537
$ws.Range("C80").Value = 42
$ws.Range("D80").Value = 342
$ws.Range("E80").Value = 1188
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 7

$ws.Range("A81").Value = "Cuba"
$ws.Range("B81").Value = 1537
$ws.Range("C81").Value = 36
$ws.Range("D81").Value = 714
$ws.Range("E81").Value = 759
$ws.Range("F81").Value = 10
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 64

# Row 153 - Aruba
$ws.Range("D153").Value = 81
$ws.Range("E153").Value = 17
